$d = $word.ActiveDocument

# 1. Replace the unique first sentence (full rewrite, not a simple word substitution).
$d.Content.Find.Execute(
    "Opera extensions -> 오페라 확장이 아닌 오페라 애드온으로 했읍니다. 이쪽이 좀 더 입에 잘 붙어서",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "extension -> 확장프로그램으로 번역했습니다.",
    2)

# 2. Replace every remaining occurrence of "애드온" with "확장프로그램"
#    throughout the rest of the document.
$d.Content.Find.Execute(
    "애드온", $true, $false, $false, $false, $false, $true, 1, $false,
    "확장프로그램",
    2)
